{"js": "// This worksheet regenerates its multiplication problems on each run;\n// every \"AA\u00d7BB=\" expression that appears is unique in the document, so\n// matching the old expression text and replacing it in place unambiguously\n// retargets the correct table cell (formatting/run structure is untouched).\nconst replacements = [\n  [\"29\u00d734=\", \"83\u00d721=\"],\n  [\"90\u00d784=\", \"89\u00d794=\"],\n  [\"86\u00d771=\", \"61\u00d737=\"],\n  [\"43\u00d771=\", \"26\u00d767=\"],\n  [\"13\u00d726=\", \"64\u00d736=\"],\n  [\"29\u00d740=\", \"65\u00d774=\"],\n  [\"92\u00d786=\", \"71\u00d777=\"],\n  [\"97\u00d716=\", \"56\u00d769=\"],\n  [\"69\u00d770=\", \"81\u00d788=\"],\n  [\"54\u00d798=\", \"63\u00d736=\"],\n  [\"77\u00d736=\", \"45\u00d781=\"],\n  [\"66\u00d797=\", \"53\u00d746=\"],\n  [\"60\u00d749=\", \"94\u00d776=\"],\n  [\"86\u00d757=\", \"22\u00d741=\"],\n  [\"45\u00d758=\", \"38\u00d772=\"],\n  [\"77\u00d791=\", \"17\u00d724=\"],\n  [\"41\u00d775=\", \"53\u00d734=\"],\n  [\"59\u00d799=\", \"20\u00d778=\"],\n  [\"67\u00d759=\", \"76\u00d757=\"],\n  [\"38\u00d721=\", \"33\u00d765=\"],\n  [\"47\u00d740=\", \"34\u00d784=\"],\n  [\"91\u00d798=\", \"57\u00d727=\"],\n  [\"82\u00d784=\", \"31\u00d729=\"],\n  [\"45\u00d746=\", \"11\u00d785=\"],\n  [\"91\u00d742=\", \"39\u00d776=\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length > 0) {\n    results.items[0].insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Each multiplication expression in the worksheet is unique, so a\n# simple Find/Replace-All per pair unambiguously retargets the right cell.\n$replacements = @(\n    @{ Find = \"29\u00d734=\"; Replace = \"83\u00d721=\" }\n    @{ Find = \"90\u00d784=\"; Replace = \"89\u00d794=\" }\n    @{ Find = \"86\u00d771=\"; Replace = \"61\u00d737=\" }\n    @{ Find = \"43\u00d771=\"; Replace = \"26\u00d767=\" }\n    @{ Find = \"13\u00d726=\"; Replace = \"64\u00d736=\" }\n    @{ Find = \"29\u00d740=\"; Replace = \"65\u00d774=\" }\n    @{ Find = \"92\u00d786=\"; Replace = \"71\u00d777=\" }\n    @{ Find = \"97\u00d716=\"; Replace = \"56\u00d769=\" }\n    @{ Find = \"69\u00d770=\"; Replace = \"81\u00d788=\" }\n    @{ Find = \"54\u00d798=\"; Replace = \"63\u00d736=\" }\n    @{ Find = \"77\u00d736=\"; Replace = \"45\u00d781=\" }\n    @{ Find = \"66\u00d797=\"; Replace = \"53\u00d746=\" }\n    @{ Find = \"60\u00d749=\"; Replace = \"94\u00d776=\" }\n    @{ Find = \"86\u00d757=\"; Replace = \"22\u00d741=\" }\n    @{ Find = \"45\u00d758=\"; Replace = \"38\u00d772=\" }\n    @{ Find = \"77\u00d791=\"; Replace = \"17\u00d724=\" }\n    @{ Find = \"41\u00d775=\"; Replace = \"53\u00d734=\" }\n    @{ Find = \"59\u00d799=\"; Replace = \"20\u00d778=\" }\n    @{ Find = \"67\u00d759=\"; Replace = \"76\u00d757=\" }\n    @{ Find = \"38\u00d721=\"; Replace = \"33\u00d765=\" }\n    @{ Find = \"47\u00d740=\"; Replace = \"34\u00d784=\" }\n    @{ Find = \"91\u00d798=\"; Replace = \"57\u00d727=\" }\n    @{ Find = \"82\u00d784=\"; Replace = \"31\u00d729=\" }\n    @{ Find = \"45\u00d746=\"; Replace = \"11\u00d785=\" }\n    @{ Find = \"91\u00d742=\"; Replace = \"39\u00d776=\" }\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($pair.Find, $false, $false, $false, $false, $false, $true, 1, $false, $pair.Replace, 2)\n}\n\nWrite-Output \"done\""}
